$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 1505.6666
$ws.Range("I52").Value = 1505.6666
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 4516.9998
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -4356.9998
$ws.Range("N52").Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 10247.917
$ws.Range("I62").Value = 9121.875
$ws.Range("J62").Value = 12500
$ws.Range("K62").Value = 9121.875
$ws.Range("L62").Value = 12500
$ws.Range("M62").Value = -8497.875
$ws.Range("N62").Value = -13748

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 10247.917
$ws.Range("I65").Value = 9121.875
$ws.Range("J65").Value = 12500
$ws.Range("K65").Value = 45609.375
$ws.Range("L65").Value = 62500
$ws.Range("M65").Value = -42489.375
$ws.Range("N65").Value = -68740

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 602.6
$ws.Range("I99").Value = 602.6
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1807.8
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -309.8000000000002
$ws.Range("N99").Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 40000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 40000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1050.6364
$ws.Range("I129").Value = 469
$ws.Range("J129").Value = 1154.5
$ws.Range("K129").Value = 1407
$ws.Range("L129").Value = 3463.5
$ws.Range("M129").Value = 3593
$ws.Range("N129").Value = -13463.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 450025000
$ws.Range("I10").Value = 450025000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 450025000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -450024830
$ws.Range("N10").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2570.9343
$ws.Range("I32").Value = 2243.0789
$ws.Range("J32").Value = 3112.6086
$ws.Range("K32").Value = 2243.0789
$ws.Range("L32").Value = 3112.6086
$ws.Range("M32").Value = -1956.0789
$ws.Range("N32").Value = -3686.6086

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2213.3462
$ws.Range("I61").Value = 1552.0476
$ws.Range("J61").Value = 4990.8
$ws.Range("K61").Value = 1552.0476
$ws.Range("L61").Value = 4990.8
$ws.Range("M61").Value = -1340.0476
$ws.Range("N61").Value = -5414.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 14493068
$ws.Range("I97").Value = 14493068
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 14493068
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -14492572
$ws.Range("N97").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1642.4736
$ws.Range("I122").Value = 1366.0769
$ws.Range("J122").Value = 2241.3333
$ws.Range("K122").Value = 4098.2307
$ws.Range("L122").Value = 6723.999899999999
$ws.Range("M122").Value = -1648.2307
$ws.Range("N122").Value = -11623.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 58500
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 58500
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 58500
$ws.Range("N128").Value = -68460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H129").Value = 42999.5
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 42999.5
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 42999.5
$ws.Range("N129").Value = -52999.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 46500
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 46500
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 46500
$ws.Range("N131").Value = -56580

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2213.3462
$ws.Range("I136").Value = 1552.0476
$ws.Range("J136").Value = 4990.8
$ws.Range("K136").Value = 4656.142800000001
$ws.Range("L136").Value = 14972.4
$ws.Range("M136").Value = -2106.142800000001
$ws.Range("N136").Value = -20072.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2189.8
$ws.Range("I20").Value = 1574.5
$ws.Range("J20").Value = 2600
$ws.Range("K20").Value = 1574.5
$ws.Range("L20").Value = 2600
$ws.Range("M20").Value = -1327.5
$ws.Range("N20").Value = -3094

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 50000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 50000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 50000
$ws.Range("N60").Value = -51198

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 11907932
$ws.Range("I105").Value = 16669693
$ws.Range("J105").Value = 3531.25
$ws.Range("K105").Value = 16669693
$ws.Range("L105").Value = 3531.25
$ws.Range("M105").Value = -16667946
$ws.Range("N105").Value = -7025.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 28666.666
$ws.Range("I123").Value = 26000
$ws.Range("J123").Value = 30000
$ws.Range("K123").Value = 26000
$ws.Range("L123").Value = 30000
$ws.Range("M123").Value = -21100
$ws.Range("N123").Value = -39800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 10000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 10000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -10884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 6450
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 6450
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 6450
$ws.Range("M45").Value = ""
$ws.Range("N45").Value = -7636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 10633.25
$ws.Range("I39").Value = 2000
$ws.Range("J39").Value = 11008.608
$ws.Range("K39").Value = 6000
$ws.Range("L39").Value = 33025.824
$ws.Range("M39").Value = -5706
$ws.Range("N39").Value = -33613.824

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3828.2766
$ws.Range("I68").Value = 4374.343
$ws.Range("J68").Value = 2235.5833
$ws.Range("K68").Value = 13123.029
$ws.Range("L68").Value = 6706.749899999999
$ws.Range("M68").Value = -12312.029
$ws.Range("N68").Value = -8328.749899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 3828.2766
$ws.Range("I71").Value = 4374.343
$ws.Range("J71").Value = 2235.5833
$ws.Range("K71").Value = 39369.087
$ws.Range("L71").Value = 20120.2497
$ws.Range("M71").Value = -35313.087
$ws.Range("N71").Value = -28232.2497

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2515.8374
$ws.Range("I131").Value = 452
$ws.Range("J131").Value = 2810.6714
$ws.Range("K131").Value = 1356
$ws.Range("L131").Value = 8432.014200000001
$ws.Range("M131").Value = 3684
$ws.Range("N131").Value = -18512.0142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 6181.727
$ws.Range("I133").Value = 4624.875
$ws.Range("J133").Value = 10333.333
$ws.Range("K133").Value = 13874.625
$ws.Range("L133").Value = 30999.999
$ws.Range("M133").Value = -8814.625
$ws.Range("N133").Value = -41119.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2223522.2
$ws.Range("I122").Value = 2779152.8
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 8337458.399999999
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -8335008.399999999
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2794.2307
$ws.Range("I7").Value = 1250
$ws.Range("J7").Value = 3257.5
$ws.Range("K7").Value = 1250
$ws.Range("L7").Value = 3257.5
$ws.Range("M7").Value = -1138
$ws.Range("N7").Value = -3481.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 9895
$ws.Range("I61").Value = 18000
$ws.Range("J61").Value = 1790
$ws.Range("K61").Value = 18000
$ws.Range("L61").Value = 1790
$ws.Range("M61").Value = -17798
$ws.Range("N61").Value = -2194

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 9895
$ws.Range("I113").Value = 18000
$ws.Range("J113").Value = 1790
$ws.Range("K113").Value = 18000
$ws.Range("L113").Value = 1790
$ws.Range("M113").Value = -15830
$ws.Range("N113").Value = -6130

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2794.2307
$ws.Range("I126").Value = 1250
$ws.Range("J126").Value = 3257.5
$ws.Range("K126").Value = 3750
$ws.Range("L126").Value = 9772.5
$ws.Range("M126").Value = -1280
$ws.Range("N126").Value = -14712.5
